$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row for the CAMP_TITLE entry, right after CAMP_RETURN_TO_PLANETE (row 2),
# pushing every following row down by one.
$null = $ws.Rows.Item(3).Insert()

$ws.Range("A3").Value = "CAMP_TITLE"
$ws.Range("B3").Value = "Welcome to the camp"
$ws.Range("C3").Value = "Bienvenu au campement"

# Fix the French translation wording for OBJECTIVE_COLLECT_FOR_TENTS (now shifted to row 27):
# "camp" -> "campement"
$ws.Range("C27").Value = "Afin de construire construire votre campement : récoltez 10 unités de pierre."

# Freeze the header row (row 1) and restore the view's selection state
$excel.ActiveWindow.DisplayGridlines = $true
$null = $ws.Range("A1").Select()
$null = $ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$null = $ws.Range("C29").Select()
